$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 438, shifting existing rows 438-467 down to 439-468.
$ws.Rows.Item(438).Insert()

# Populate the newly inserted row 438 with the new data record.
$ws.Cells.Item(438, 1).Value = 3
$ws.Cells.Item(438, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(438, 3).Value = "Coquimbo"
$ws.Cells.Item(438, 4).Value = 44931
$ws.Cells.Item(438, 5).Value = 5
$ws.Cells.Item(438, 6).Value = 100112009
$ws.Cells.Item(438, 7).Value = "Acelga"
$ws.Cells.Item(438, 8).Value = "Sin especificar"
$ws.Cells.Item(438, 9).Value = "Primera"
$ws.Cells.Item(438, 10).Value = 210
$ws.Cells.Item(438, 11).Value = 4000
$ws.Cells.Item(438, 12).Value = 4500
$ws.Cells.Item(438, 13).Value = 4262
$ws.Cells.Item(438, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(438, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(438, 16).Value = 710
$ws.Cells.Item(438, 17).Value = 6
$ws.Cells.Item(438, 18).Value = "Hortaliza"
